$d = $word.ActiveDocument

# Replaces the run content of an EXISTING paragraph (keeps its w:p attributes
# and w:pPr untouched) by targeting the range up to, but excluding, the
# paragraph mark.
function Set-ParagraphRuns($para, $runsXml) {
    $full = $para.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $runsXml + '</w:p>'
    $r.InsertXML($xml) | Out-Null
}

# Fills a brand-new (empty) paragraph: targets the paragraph's full range
# (including its mark) so no stray empty run is left behind, and supplies an
# explicit w:pPr so the paragraph-mark run properties are preserved too.
function Set-NewParagraphRuns($para, $runsXml) {
    $full = $para.Range
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + $runsXml + '</w:p>'
    $full.InsertXML($xml) | Out-Null
}

# --- Paragraph 1: "Dear Mr X," -> "Dear " + "Mr" (spell-checked) + " X," ---
$p1Runs = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Dear </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Mr</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> X,</w:t></w:r>'
Set-ParagraphRuns $d.Paragraphs(1) $p1Runs

# --- Paragraph 2: "{{ hook_paragraph }}" -> "{{ " + "opening_paragraph" (spell-checked) + " " + "}}" ---
$p2Runs = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>opening_paragraph</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r>'
Set-ParagraphRuns $d.Paragraphs(2) $p2Runs

# --- Paragraph 3: "{{ main_section }}" -> "{{ " + "core"+"_"+"parag"+"r"+"aphs" (spell-checked) + " }}" ---
$p3Runs = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>core</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>parag</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>aphs</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>'
Set-ParagraphRuns $d.Paragraphs(3) $p3Runs

# --- New paragraph 4 (inserted after paragraph 3, before "Yours sincerely,"):
#     "{{ closing_paragraph }}" -> "{{ " + "c"+"losing"+"_par"+"a"+"g"+"r"+"aph" (spell-checked) + " }}" ---
$d.Paragraphs(3).Range.InsertParagraphAfter() | Out-Null
$p4Runs = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>c</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>losing</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>_par</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>a</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>g</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>aph</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>'
Set-NewParagraphRuns $d.Paragraphs(4) $p4Runs
